# Updated TODO list after meeting
#
# The most recent week's sheet ("12-3-14 to 12-13-14") is duplicated to
# create a new current-week sheet ("12-12-14 to 12-19-14") placed in front
# of it. The new sheet keeps only the section headers plus a handful of
# freshly defined tasks (the per-person/definition-of-finished detail for
# those tasks has not been filled in yet). The old sheet is renamed to
# reflect that its period now ends on 12-12-14 instead of 12-13-14, and its
# content is left untouched.

$wb = $excel.ActiveWorkbook

$mostRecent = $wb.Worksheets.Item(1)

# Duplicate the most recent week's sheet; the copy is placed immediately
# before it and becomes the new active/first sheet.
$mostRecent.Copy($mostRecent)

$newWeek = $wb.Worksheets.Item(1)
$oldWeek = $wb.Worksheets.Item(2)

$newWeek.Name = "12-12-14 to 12-19-14"
$oldWeek.Name = "12-3-14 to 12-12-14"

# Clear out the detailed rows in the new week's "Design Tasks" section
# (rows 3-10) and replace them with the newly identified tasks, which only
# have a task name filled in so far (no assignee or definition of finished
# yet).
$newWeek.Range("B3:C10").Clear()
$newWeek.Range("A3").Value = "Derive 27 equations"
$newWeek.Range("A4").Value = "Estimate component masses"
$newWeek.Range("A5").Value = "Shoulder joint design"
$newWeek.Range("A6").Value = "Stress test the microcontroller"
$newWeek.Range("A7").Value = "Look into master slave architecture"

# Rows 8-11 (the rest of the old Design Tasks detail) are no longer used.
$newWeek.Range("A8:C11").Clear()

# The last row of the Report Tasks section (old row 16) is also removed.
$newWeek.Range("A16:C16").Clear()

# Match the author's last on-sheet selection.
$newWeek.Activate()
$null = $newWeek.Range("A16").Select()
